# Add the new intro paragraphs + "11 I Need to Escape" heading right after
# the "Write Up" title paragraph, leaving the trailing empty paragraph and
# the trailing empty Heading1 paragraph untouched at the end of the body.

$d = $word.ActiveDocument

# --- Create the three new paragraphs (chained off the Title paragraph) ---
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$d.Paragraphs(2).Style = "Normal"
$introPara = $d.Paragraphs(2)
$introPara.Range.InsertParagraphAfter() | Out-Null

$d.Paragraphs(3).Style = "Normal"
$joinUsPara = $d.Paragraphs(3)
$joinUsPara.Range.InsertParagraphAfter() | Out-Null

$d.Paragraphs(4).Style = "Heading1"

# --- Paragraph 2: the "In this section..." body text ---
$introPara = $d.Paragraphs(2)
$introRange = $introPara.Range
$introRange.Collapse(1)
$introRange.InsertAfter("In this section we will be creating")
$introRange.Collapse(0)
$introRange.InsertAfter(" a bit of error fixing code, in case our hero gets into trouble, and finds that he spawn")
$introRange.Collapse(0)
$introRange.InsertAfter("ed")
$introRange.Collapse(0)
$introRange.InsertAfter(" himself right inside of the wall inside of a dungeon.")
$introRange.Collapse(0)
$introRange.InsertAfter(" This code, when initiated, will cause the hero to jump out of the wall and into an open corridor.")

# --- Paragraph 3: the "So, if this is something..." text ---
$joinUsPara = $d.Paragraphs(3)
$joinUsRange = $joinUsPara.Range
$joinUsRange.Collapse(1)
$joinUsRange.InsertAfter("So, if this is something that you would like to learn just a little bit more about, then please join us for our brand-new article, entitled:")

# --- Paragraph 4: the new Heading1 "11 I Need to Escape" ---
$headingPara = $d.Paragraphs(4)
$headingRange = $headingPara.Range
$headingRange.Collapse(1)
$headingRange.InsertAfter("11 I Need to Escape")

Write-Output "Inserted intro paragraphs and new heading."
